$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 993.5294
$ws.Range("I19").Value = 639.8
$ws.Range("J19").Value = 1140.9166
$ws.Range("K19").Value = 639.8
$ws.Range("L19").Value = 1140.9166
$ws.Range("M19").Value = -464.8
$ws.Range("N19").Value = -1490.9166
$ws.Range("H40").Value = 2115.889
$ws.Range("I40").Value = 2037.7693
$ws.Range("J40").Value = 2319
$ws.Range("K40").Value = 2037.7693
$ws.Range("L40").Value = 2319
$ws.Range("M40").Value = -1862.7693
$ws.Range("N40").Value = -2669
$ws.Range("H116").Value = 1566.6666
$ws.Range("I116").Value = 1600
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 1600
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1842
$ws.Range("N116").Value = -8384
$ws.Range("H137").Value = 1350.4762
$ws.Range("I137").Value = 1088.421
$ws.Range("J137").Value = 3840
$ws.Range("K137").Value = 3265.263
$ws.Range("L137").Value = 11520
$ws.Range("M137").Value = -715.2629999999999
$ws.Range("N137").Value = -16620
$ws.Range("H138").Value = 3611.2876
$ws.Range("I138").Value = 1289.0968
$ws.Range("J138").Value = 5325.2856
$ws.Range("K138").Value = 3867.2904
$ws.Range("L138").Value = 15975.8568
$ws.Range("M138").Value = 1272.7096
$ws.Range("N138").Value = -26255.8568
$ws.Range("H141").Value = 2977.3809
$ws.Range("I141").Value = 2796.0527
$ws.Range("J141").Value = 4700
$ws.Range("K141").Value = 8388.158100000001
$ws.Range("L141").Value = 14100
$ws.Range("M141").Value = -3208.158100000001
$ws.Range("N141").Value = -24460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 975
$ws.Range("I21").Value = 975
$ws.Range("K21").Value = 975
$ws.Range("M21").Value = -601
$ws.Range("H32").Value = 23307.75
$ws.Range("I32").Value = 4052.0747
$ws.Range("J32").Value = 166655.56
$ws.Range("K32").Value = 4052.0747
$ws.Range("L32").Value = 166655.56
$ws.Range("M32").Value = -3765.0747
$ws.Range("N32").Value = -167229.56
$ws.Range("H45").Value = 2078.7
$ws.Range("I45").Value = 1462.5
$ws.Range("J45").Value = 2489.5
$ws.Range("K45").Value = 1462.5
$ws.Range("L45").Value = 2489.5
$ws.Range("M45").Value = -1085.5
$ws.Range("N45").Value = -3243.5
$ws.Range("H61").Value = 3038.5454
$ws.Range("I61").Value = 2901.8333
$ws.Range("J61").Value = 3202.6
$ws.Range("K61").Value = 2901.8333
$ws.Range("L61").Value = 3202.6
$ws.Range("M61").Value = -2689.8333
$ws.Range("N61").Value = -3626.6
$ws.Range("H74").Value = 1678.5862
$ws.Range("I74").Value = 1159.0952
$ws.Range("J74").Value = 3042.25
$ws.Range("K74").Value = 1159.0952
$ws.Range("L74").Value = 3042.25
$ws.Range("M74").Value = -285.0952
$ws.Range("N74").Value = -4790.25
$ws.Range("H77").Value = 1678.5862
$ws.Range("I77").Value = 1159.0952
$ws.Range("J77").Value = 3042.25
$ws.Range("K77").Value = 5795.476
$ws.Range("L77").Value = 15211.25
$ws.Range("M77").Value = -1427.476
$ws.Range("N77").Value = -23947.25
$ws.Range("H132").Value = 1845.1562
$ws.Range("I132").Value = 1451.24
$ws.Range("K132").Value = 4353.72
$ws.Range("M132").Value = -1823.72
$ws.Range("H136").Value = 3038.5454
$ws.Range("I136").Value = 2901.8333
$ws.Range("J136").Value = 3202.6
$ws.Range("K136").Value = 8705.499899999999
$ws.Range("L136").Value = 9607.799999999999
$ws.Range("M136").Value = -6155.499899999999
$ws.Range("N136").Value = -14707.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 79052.30499999999
$ws.Range("I20").Value = 127337.5
$ws.Range("J20").Value = 1796
$ws.Range("K20").Value = 127337.5
$ws.Range("L20").Value = 1796
$ws.Range("M20").Value = -127090.5
$ws.Range("N20").Value = -2290
$ws.Range("H107").Value = 62501052
$ws.Range("I107").Value = 111112296
$ws.Range("K107").Value = 111112296
$ws.Range("M107").Value = -111110376
$ws.Range("H108").Value = 35999
$ws.Range("J108").Value = 35999
$ws.Range("L108").Value = 35999
$ws.Range("N108").Value = -43679
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2010
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 2010
$ws.Range("K6").Value = 0
$ws.Range("L6").ClearContents()
$ws.Range("M6").Value = 2010
$ws.Range("N6").Value = -2236
$ws.Range("H31").Value = 33731.938
$ws.Range("I31").Value = 56933.055
$ws.Range("J31").Value = 3901.9285
$ws.Range("K31").Value = 56933.055
$ws.Range("L31").Value = 3901.9285
$ws.Range("M31").Value = -56638.055
$ws.Range("N31").Value = -4491.9285
$ws.Range("H34").Value = 33731.938
$ws.Range("I34").Value = 56933.055
$ws.Range("J34").Value = 3901.9285
$ws.Range("K34").Value = 56933.055
$ws.Range("L34").Value = 3901.9285
$ws.Range("M34").Value = -56731.055
$ws.Range("N34").Value = -4305.9285
$ws.Range("H58").Value = 12651.739
$ws.Range("I58").Value = 1482.9474
$ws.Range("J58").Value = 65703.5
$ws.Range("K58").Value = 1482.9474
$ws.Range("L58").Value = 65703.5
$ws.Range("M58").Value = -1279.9474
$ws.Range("N58").Value = -66109.5
$ws.Range("H99").Value = 9030.933999999999
$ws.Range("I99").Value = 2554
$ws.Range("J99").Value = 12269.4
$ws.Range("K99").Value = 2554
$ws.Range("L99").Value = 12269.4
$ws.Range("M99").Value = -1056
$ws.Range("N99").Value = -15265.4
$ws.Range("H126").Value = 9030.933999999999
$ws.Range("I126").Value = 2554
$ws.Range("J126").Value = 12269.4
$ws.Range("K126").Value = 7662
$ws.Range("L126").Value = 36808.2
$ws.Range("M126").Value = -5192
$ws.Range("N126").Value = -41748.2
$ws.Range("H132").Value = 62504340
$ws.Range("I132").Value = 76928700
$ws.Range("J132").Value = 45457370
$ws.Range("K132").Value = 230786100
$ws.Range("L132").Value = 136372110
$ws.Range("M132").Value = -230783570
$ws.Range("N132").Value = -136377170
$ws.Range("H134").Value = 1496.5416
$ws.Range("I134").Value = 1450.7727
$ws.Range("K134").Value = 4352.3181
$ws.Range("M134").Value = -1817.3181
$ws.Range("H136").Value = 12651.739
$ws.Range("I136").Value = 1482.9474
$ws.Range("J136").Value = 65703.5
$ws.Range("K136").Value = 4448.8422
$ws.Range("L136").Value = 197110.5
$ws.Range("M136").Value = -1898.8422
$ws.Range("N136").Value = -202210.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6831.1143
$ws.Range("I5").Value = 1078.5
$ws.Range("J5").Value = 19382.273
$ws.Range("K5").Value = 3235.5
$ws.Range("L5").Value = 58146.819
$ws.Range("M5").Value = -3123.5
$ws.Range("N5").Value = -58370.819
$ws.Range("H36").Value = 325.25
$ws.Range("I36").Value = 267
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 801
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -632
$ws.Range("N36").Value = -1838
$ws.Range("H69").Value = 2666.6667
$ws.Range("J69").Value = 2666.6667
$ws.Range("L69").Value = 8000.000100000001
$ws.Range("N69").Value = -9622.000100000001
$ws.Range("H72").Value = 2666.6667
$ws.Range("J72").Value = 2666.6667
$ws.Range("L72").Value = 24000.0003
$ws.Range("N72").Value = -32112.0003
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H135").Value = 6831.1143
$ws.Range("I135").Value = 1078.5
$ws.Range("J135").Value = 19382.273
$ws.Range("K135").Value = 9706.5
$ws.Range("L135").Value = 174440.457
$ws.Range("M135").Value = -7171.5
$ws.Range("N135").Value = -179510.457

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2777.4
$ws.Range("I132").Value = 2151.2778
$ws.Range("J132").Value = 4387.4287
$ws.Range("K132").Value = 6453.8334
$ws.Range("L132").Value = 13162.2861
$ws.Range("M132").Value = -3923.8334
$ws.Range("N132").Value = -18222.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2576.1177
$ws.Range("I7").Value = 1929.9
$ws.Range("J7").Value = 3499.2856
$ws.Range("K7").Value = 1929.9
$ws.Range("L7").Value = 3499.2856
$ws.Range("M7").Value = -1817.9
$ws.Range("N7").Value = -3723.2856
$ws.Range("H16").Value = 63238.062
$ws.Range("I16").Value = 143557.42
$ws.Range("J16").Value = 767.44446
$ws.Range("K16").Value = 143557.42
$ws.Range("L16").Value = 767.44446
$ws.Range("M16").Value = -143387.42
$ws.Range("N16").Value = -1107.44446
$ws.Range("H33").Value = 15555
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 15555
$ws.Range("K33").Value = 0
$ws.Range("L33").ClearContents()
$ws.Range("M33").Value = 15555
$ws.Range("N33").Value = -16135
$ws.Range("H40").Value = 102209.9
$ws.Range("I40").Value = 501000
$ws.Range("J40").Value = 2512.375
$ws.Range("K40").Value = 501000
$ws.Range("L40").Value = 2512.375
$ws.Range("M40").Value = -500864
$ws.Range("N40").Value = -2784.375
$ws.Range("H69").Value = 40000
$ws.Range("I69").Value = 40000
$ws.Range("K69").Value = 40000
$ws.Range("M69").Value = -39189
$ws.Range("H72").Value = 40000
$ws.Range("I72").Value = 40000
$ws.Range("K72").Value = 120000
$ws.Range("M72").Value = -115944
$ws.Range("H126").Value = 2576.1177
$ws.Range("I126").Value = 1929.9
$ws.Range("J126").Value = 3499.2856
$ws.Range("K126").Value = 5789.700000000001
$ws.Range("L126").Value = 10497.8568
$ws.Range("M126").Value = -3319.700000000001
$ws.Range("N126").Value = -15437.8568
$ws.Range("H132").Value = 4730.7334
$ws.Range("I132").Value = 5450.8
$ws.Range("J132").Value = 3290.6
$ws.Range("K132").Value = 16352.4
$ws.Range("L132").Value = 9871.799999999999
$ws.Range("M132").Value = -13822.4
$ws.Range("N132").Value = -14931.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3750.625
$ws.Range("J122").Value = 7001.6665
$ws.Range("L122").Value = 21004.9995
$ws.Range("N122").Value = -25904.9995
$ws.Range("H132").Value = 4661.8125
$ws.Range("I132").Value = 5787.8887
$ws.Range("J132").Value = 3214
$ws.Range("K132").Value = 17363.6661
$ws.Range("L132").Value = 9642
$ws.Range("M132").Value = -14833.6661
$ws.Range("N132").Value = -14702
$ws.Range("H136").Value = 1161.5745
$ws.Range("I136").Value = 464.32257
$ws.Range("J136").Value = 2512.5
$ws.Range("K136").Value = 2512.5
$ws.Range("L136").Value = 7537.5
$ws.Range("M136").Value = 1157.03229
$ws.Range("N136").Value = -12637.5
